$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# and mangles the exact printed representation (trailing zeros, float noise).
$textCells = @("D4","D5","D6","D9","D11","D12","D13","D16","D18","D19","D21","D22","D23","D24","D25","D28","D29","D30","D31","D33","D34","D36","D37","D38","D40","D43","D49","D50","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '69.399.55'
$ws.Range("E2").Value = '  +1.73%  '
$ws.Range("D3").Value = '3.946.16'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '492.07'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").Value = '147.07'
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.737'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("E10").Value = '  +4.23%  '
$ws.Range("D11").Value = '0.0000349'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("D12").Value = '43.15'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").Value = '10.47'
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("D14").Value = '4.572.58'
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").Value = '3.927.42'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = '14.29'
$ws.Range("E16").Value = '  -4.09%  '
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").Value = '19.95'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '1.17'
$ws.Range("E19").Value = '  +2.21%  '
$ws.Range("D20").Value = '69.442.65'
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").Value = '441.04'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = '3.47'
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("D23").Value = '14.52'
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("D24").Value = '89.40'
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("D25").Value = '11.99'
$ws.Range("E25").Value = '  +8.81%  '
$ws.Range("E26").Value = '  +2.79%  '
$ws.Range("E27").Value = '  -4.29%  '
$ws.Range("D28").Value = '37.28'
$ws.Range("E28").Value = '  -4.67%  '
$ws.Range("D29").Value = '5.65'
$ws.Range("E29").Value = '  -3.75%  '
$ws.Range("D30").Value = '707.80'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '13.51'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '2.90'
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").Value = '0.474'
$ws.Range("E34").Value = '  +26.34%  '
$ws.Range("D35").Value = '0.0₃0916'
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("D36").Value = '61.55'
$ws.Range("E36").Value = '  +4.41%  '
$ws.Range("D37").Value = '6.06'
$ws.Range("E37").Value = '  +4.30%  '
$ws.Range("D38").Value = '40.89'
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("E42").Value = '  +2.41%  '
$ws.Range("D43").Value = '2.94'
$ws.Range("E43").Value = '  +2.41%  '
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("E45").Value = '  +2.02%  '
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '0.0₆0366'
$ws.Range("E47").Value = '  +9.23%  '
$ws.Range("E48").Value = '  +6.42%  '
$ws.Range("D49").Value = '3.07'
$ws.Range("E49").Value = '  +7.84%  '
$ws.Range("D50").Value = '3.39'
$ws.Range("E50").Value = '  -1.16%  '
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").Value = '2.08'
$ws.Range("E51").Value = '  -3.25%  '
